$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $text)
    # Temporarily force Text format so the numeric-looking string is stored
    # as a text value (matching the workbook's inlineStr cells), then restore
    # the cell style to Normal so no stray number-format style is left behind.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Row 2 (IIT VAL)
Set-TextValue $ws.Range("B2") "0.009173228"
Set-TextValue $ws.Range("C2") "0.078686036"

# Row 3 (REG VAL)
Set-TextValue $ws.Range("B3") "0.014684563"
Set-TextValue $ws.Range("C3") "0.0831223"
$ws.Range("D3").Value = 77
$ws.Range("E3").Value = 3

# Row 4 (IIT TEST)
Set-TextValue $ws.Range("B4") "0.021424113"
Set-TextValue $ws.Range("C4") "0.110698275"

# Row 5 (REG TEST)
Set-TextValue $ws.Range("B5") "0.01987446"
Set-TextValue $ws.Range("C5") "0.11125524"
